$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 832.5122
$ws.Range("I6").Value = 252.71428
$ws.Range("J6").Value = 951.8823
$ws.Range("K6").Value = 758.14284
$ws.Range("L6").Value = 2855.6469
$ws.Range("M6").Value = -646.14284
$ws.Range("N6").Value = -3079.6469

$ws.Range("H112").Value = 3544.5417
$ws.Range("J112").Value = 3644.3914
$ws.Range("L112").Value = 10933.1742
$ws.Range("N112").Value = -13149.1742

$ws.Range("H137").Value = 2824
$ws.Range("I137").Value = 2430.6667
$ws.Range("J137").Value = 3217.3333
$ws.Range("K137").Value = 7292.000100000001
$ws.Range("L137").Value = 9651.999899999999
$ws.Range("M137").Value = -4742.000100000001
$ws.Range("N137").Value = -14751.9999

$ws.Range("H138").Value = 7861.711
$ws.Range("J138").Value = 9301.529
$ws.Range("L138").Value = 27904.587
$ws.Range("N138").Value = -38184.587

$ws.Range("H141").Value = 4340.8184
$ws.Range("I141").Value = 4456.875
$ws.Range("J141").Value = 4031.3333
$ws.Range("K141").Value = 13370.625
$ws.Range("L141").Value = 12093.9999
$ws.Range("M141").Value = -8190.625
$ws.Range("N141").Value = -22453.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6028.4463
$ws.Range("I32").Value = 641.7593000000001
$ws.Range("J32").Value = 32472.182
$ws.Range("K32").Value = 641.7593000000001
$ws.Range("L32").Value = 32472.182
$ws.Range("M32").Value = -354.7593000000001
$ws.Range("N32").Value = -33046.182

$ws.Range("H61").Value = 3250.5334
$ws.Range("I61").Value = 3142.5908
$ws.Range("K61").Value = 3142.5908
$ws.Range("M61").Value = -2930.5908

$ws.Range("H97").Value = 2998.842
$ws.Range("I97").Value = 2153.125
$ws.Range("J97").Value = 3613.9092
$ws.Range("K97").Value = 2153.125
$ws.Range("L97").Value = 3613.9092
$ws.Range("M97").Value = -1657.125
$ws.Range("N97").Value = -4605.9092

$ws.Range("H122").Value = 2361
$ws.Range("J122").Value = 2361
$ws.Range("L122").Value = 7083
$ws.Range("N122").Value = -11983

$ws.Range("H136").Value = 3250.5334
$ws.Range("I136").Value = 3142.5908
$ws.Range("K136").Value = 9427.7724
$ws.Range("M136").Value = -6877.7724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1198.1892
$ws.Range("I20").Value = 1328.68
$ws.Range("J20").Value = 926.3333
$ws.Range("K20").Value = 1328.68
$ws.Range("L20").Value = 926.3333
$ws.Range("M20").Value = -1081.68
$ws.Range("N20").Value = -1420.3333

$ws.Range("H80").Value = 1093.08
$ws.Range("J80").Value = 645.75
$ws.Range("L80").Value = 645.75
$ws.Range("N80").Value = -2641.75

$ws.Range("H83").Value = 1093.08
$ws.Range("J83").Value = 645.75
$ws.Range("L83").Value = 3228.75
$ws.Range("N83").Value = -13212.75

$ws.Range("H86").Value = 3053.3635
$ws.Range("I86").Value = 2202.625
$ws.Range("K86").Value = 2202.625
$ws.Range("M86").Value = -1079.625

$ws.Range("H89").Value = 3053.3635
$ws.Range("I89").Value = 2202.625
$ws.Range("K89").Value = 11013.125
$ws.Range("M89").Value = -5397.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5857.2666
$ws.Range("I31").Value = 4863.4287
$ws.Range("J31").Value = 8176.222
$ws.Range("K31").Value = 4863.4287
$ws.Range("L31").Value = 8176.222
$ws.Range("M31").Value = -4568.4287
$ws.Range("N31").Value = -8766.222

$ws.Range("H34").Value = 5857.2666
$ws.Range("I34").Value = 4863.4287
$ws.Range("J34").Value = 8176.222
$ws.Range("K34").Value = 4863.4287
$ws.Range("L34").Value = 8176.222
$ws.Range("M34").Value = -4661.4287
$ws.Range("N34").Value = -8580.222

$ws.Range("H122").Value = 3422.8823
$ws.Range("I122").Value = 3242.0667
$ws.Range("J122").Value = 4779
$ws.Range("K122").Value = 9726.2001
$ws.Range("L122").Value = 14337
$ws.Range("M122").Value = -7276.2001
$ws.Range("N122").Value = -19237

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 100041800
$ws.Range("J37").Value = 100041800
$ws.Range("L37").Value = 300125400
$ws.Range("N37").Value = -300125624

$ws.Range("H68").Value = 3115.5293
$ws.Range("J68").Value = 3156.2666
$ws.Range("L68").Value = 9468.799800000001
$ws.Range("N68").Value = -11090.7998

$ws.Range("H71").Value = 3115.5293
$ws.Range("J71").Value = 3156.2666
$ws.Range("L71").Value = 28406.3994
$ws.Range("N71").Value = -36518.39939999999

$ws.Range("H121").Value = 71958.64
$ws.Range("I121").Value = 286
$ws.Range("K121").Value = 858
$ws.Range("M121").Value = 452

$ws.Range("H134").Value = 5394.1333
$ws.Range("I134").Value = 1752.75
$ws.Range("J134").Value = 19959.666
$ws.Range("K134").Value = 5258.25
$ws.Range("L134").Value = 59878.99800000001
$ws.Range("M134").Value = -188.25
$ws.Range("N134").Value = -70018.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2106.1765
$ws.Range("I102").Value = 2150.4285
$ws.Range("J102").Value = 1899.6666
$ws.Range("K102").Value = 2150.4285
$ws.Range("L102").Value = 1899.6666
$ws.Range("M102").Value = -528.4285
$ws.Range("N102").Value = -5143.6666

$ws.Range("H132").Value = 3437.5134
$ws.Range("I132").Value = 3567.724
$ws.Range("J132").Value = 2965.5
$ws.Range("K132").Value = 10703.172
$ws.Range("L132").Value = 8896.5
$ws.Range("M132").Value = -8173.172
$ws.Range("N132").Value = -13956.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 32000
$ws.Range("J5").Value = 32000
$ws.Range("L5").Value = 32000
$ws.Range("N5").Value = -32226

$ws.Range("H7").Value = 18798.059
$ws.Range("I7").Value = 13658.8
$ws.Range("J7").Value = 26139.857
$ws.Range("K7").Value = 13658.8
$ws.Range("L7").Value = 26139.857
$ws.Range("M7").Value = -13546.8
$ws.Range("N7").Value = -26363.857

$ws.Range("H20").Value = 1923235.2
$ws.Range("I20").Value = 1923235.2
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1923235.2
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1923009.2
$ws.Range("N20").ClearContents()

$ws.Range("H40").Value = 7133.0884
$ws.Range("I40").Value = 3838.476
$ws.Range("J40").Value = 12455.154
$ws.Range("K40").Value = 3838.476
$ws.Range("L40").Value = 12455.154
$ws.Range("M40").Value = -3702.476
$ws.Range("N40").Value = -12727.154

$ws.Range("H46").Value = 7656.125
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 7312.25
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 7312.25
$ws.Range("M46").Value = -7812
$ws.Range("N46").Value = -7688.25

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H126").Value = 18798.059
$ws.Range("I126").Value = 13658.8
$ws.Range("J126").Value = 26139.857
$ws.Range("K126").Value = 40976.39999999999
$ws.Range("L126").Value = 78419.571
$ws.Range("M126").Value = -38506.39999999999
$ws.Range("N126").Value = -83359.571

$ws.Range("H136").Value = 4375.8
$ws.Range("I136").Value = 4266.7
$ws.Range("J136").Value = 4594
$ws.Range("K136").Value = 12800.1
$ws.Range("L136").Value = 13782
$ws.Range("M136").Value = -10250.1
$ws.Range("N136").Value = -18882

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 3909
$ws.Range("I30").Value = 3909
$ws.Range("K30").Value = 3909
$ws.Range("M30").Value = -3802

$ws.Range("H122").Value = 3415.9473
$ws.Range("I122").Value = 2548.6667
$ws.Range("K122").Value = 7646.000100000001
$ws.Range("M122").Value = -5196.000100000001

$ws.Range("H132").Value = 3202.8215
$ws.Range("I132").Value = 2903.2917
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8709.875100000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6179.875100000001
$ws.Range("N132").Value = -20060
